$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet (SVM) and name it "Sheet3"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Sheet3"

# Populate the table: Metode / Akurasi header row, then KNN and SVM rows.
# Order of assignment matters for shared-string allocation order (Metode, KNN, SVM, Akurasi).
$ws3.Range("A1").Value = "Metode"
$ws3.Range("A2").Value = "KNN"
$ws3.Range("A3").Value = "SVM"
$ws3.Range("B1").Value = "Akurasi"
$ws3.Range("B2").Value = 79.411764705882348
$ws3.Range("B3").Value = 88.235294117647058

# Set column widths as close as possible to the authored widths (7.453125 / 11.453125 chars)
$ws3.Columns.Item(1).ColumnWidth = 6.71
$ws3.Columns.Item(2).ColumnWidth = 10.71
